# Swap the contents of columns C (codeforiati:group-code) and D (codeforiati:group-name)
# for every row in the used range, so that the two columns (header + all values)
# exchange places - matching the reordering applied to the shared strings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value2
    $dValue = $dCell.Value2

    $cCell.Value2 = $dValue
    $dCell.Value2 = $cValue
}
